# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.506.40"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3
$ws.Range("D3").Value = "1.649.46"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.52%  "

# Row 5
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.50%  "

# Row 6
$ws.Range("D6").Value = "'300.24"
$ws.Range("E6").Value = "  -1.20%  "

# Row 7
$ws.Range("D7").Value = "'0.3794"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3571"
$ws.Range("E8").Value = "  -1.02%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'50.68"
$ws.Range("E9").Value = "  -2.45%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.08107"
$ws.Range("E10").Value = "  -0.99%  "

# Row 11
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'1.226"
$ws.Range("E11").Value = "  -1.53%  "

# Row 12
$ws.Range("E12").Value = "  +0.39%  "

# Row 13
$ws.Range("D13").Value = "'22.09"
$ws.Range("E13").Value = "  -1.75%  "

# Row 14
$ws.Range("D14").Value = "'6.428"
$ws.Range("E14").Value = "  -1.61%  "

# Row 15
$ws.Range("D15").Value = "'7.431"
$ws.Range("E15").Value = "  +0.88%  "

# Row 16
$ws.Range("D16").Value = "'0.00001205"
$ws.Range("E16").Value = "  -2.03%  "

# Row 17
$ws.Range("D17").Value = "1.659.02"
$ws.Range("E17").Value = "  +0.88%  "

# Row 18
$ws.Range("D18").Value = "'97.29"
$ws.Range("E18").Value = "  +0.53%  "

# Row 19
$ws.Range("D19").Value = "'0.06994"
$ws.Range("E19").Value = "  +0.49%  "

# Row 20
$ws.Range("D20").Value = "'6.787"
$ws.Range("E20").Value = "  +0.63%  "

# Row 21
$ws.Range("E21").Value = "  -0.48%  "

# Row 22
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
$ws.Range("D23").Value = "'12.57"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("D24").Value = "23.532.03"
$ws.Range("E24").Value = "  -0.41%  "

# Row 25
$ws.Range("D25").Value = "'2.481"
$ws.Range("E25").Value = "  -1.52%  "

# Row 26
$ws.Range("D26").Value = "'2.925"
$ws.Range("E26").Value = "  -6.36%  "

# Row 27
$ws.Range("D27").Value = "'21.02"
$ws.Range("E27").Value = "  -1.43%  "

# Row 28
$ws.Range("D28").Value = "'152.62"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29
$ws.Range("D29").Value = "'5.232"
$ws.Range("E29").Value = "  +0.80%  "

# Row 30
$ws.Range("D30").Value = "'133.37"
$ws.Range("E30").Value = "  -1.14%  "

# Row 31
$ws.Range("D31").Value = "1.847.32"
$ws.Range("E31").Value = "  +1.01%  "

# Row 32
$ws.Range("D32").Value = "'6.971"
$ws.Range("E32").Value = "  +3.05%  "

# Row 33
$ws.Range("D33").Value = "'2.137"
$ws.Range("E33").Value = "  +4.64%  "

# Row 34
$ws.Range("D34").Value = "'11.86"
$ws.Range("E34").Value = "  +4.04%  "

# Row 35
$ws.Range("D35").Value = "'1.033"
$ws.Range("E35").Value = "  -5.83%  "

# Row 36
$ws.Range("D36").Value = "'0.02736"
$ws.Range("E36").Value = "  -1.99%  "

# Row 37
$ws.Range("D37").Value = "'0.08723"
$ws.Range("E37").Value = "  -0.81%  "

# Row 38
$ws.Range("D38").Value = "'6.000"
$ws.Range("E38").Value = "  -1.26%  "

# Row 39
$ws.Range("D39").Value = "'0.2455"
$ws.Range("E39").Value = "  -2.24%  "

# Row 40
$ws.Range("D40").Value = "'13.30"
$ws.Range("E40").Value = "  +3.69%  "

# Row 41
$ws.Range("D41").Value = "'0.06883"
$ws.Range("E41").Value = "  -1.94%  "

# Row 42
$ws.Range("D42").Value = "'0.6930"
$ws.Range("E42").Value = "  -1.91%  "

# Row 43
$ws.Range("D43").Value = "'1.323"
$ws.Range("E43").Value = "  -0.02%  "

# Row 44
$ws.Range("D44").Value = "'15.72"
$ws.Range("E44").Value = "  -0.94%  "

# Row 45
$ws.Range("D45").Value = "'0.6453"
$ws.Range("E45").Value = "  -0.95%  "

# Row 46
$ws.Range("E46").Value = "  +0.46%  "

# Row 47
$ws.Range("D47").Value = "'2.274"
$ws.Range("E47").Value = "  -2.86%  "

# Row 48
$ws.Range("D48").Value = "'3.928"
$ws.Range("E48").Value = "  -1.17%  "

# Row 49
$ws.Range("D49").Value = "'0.07821"
$ws.Range("E49").Value = "  -2.10%  "

# Row 50
$ws.Range("D50").Value = "'128.29"
$ws.Range("E50").Value = "  +0.44%  "

# Row 51
$ws.Range("D51").Value = "'1.175"
$ws.Range("E51").Value = "  -0.99%  "
